$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values in existing rows (2-21)
$ws.Range("B2").Value = "NSE:AMJLAND"
$ws.Range("C2").Value = "NSE:ACC"
$ws.Range("D2").Value = "NSE:GODREJCP"
$ws.Range("E2").Value = "NSE:ABB"
$ws.Range("F2").Value = "NSE:EXIDEIND"
$ws.Range("B3").Value = "NSE:EXIDEIND"
$ws.Range("C3").Value = "NSE:AMNPLST"
$ws.Range("D3").Value = "NSE:HDFCAMC"
$ws.Range("E3").Value = "NSE:RAMCOCEM"
$ws.Range("F3").Value = ""
$ws.Range("B4").Value = "NSE:OBCL"
$ws.Range("C4").Value = "NSE:ANDHRAPAP"
$ws.Range("D4").Value = "NSE:JUBLFOOD"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("B5").Value = "NSE:RUSHIL"
$ws.Range("C5").Value = "NSE:AWL"
$ws.Range("D5").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "NSE:BECTORFOOD"
$ws.Range("D6").Value = ""
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "NSE:CARBORUNIV"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "NSE:CASTROLIND"
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "NSE:CELEBRITY"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "NSE:CLEDUCATE"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "NSE:DATAMATICS"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "NSE:DEEPINDS"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "NSE:DELTACORP"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = "NSE:DODLA"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "NSE:EXXARO"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "NSE:FAZE3Q"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "NSE:GENCON"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "NSE:GEPIL"
$ws.Range("C19").Value = "NSE:GMBREW"
$ws.Range("C20").Value = "NSE:GNA"
$ws.Range("C21").Value = "NSE:IFCI"

# Add new rows 22-47
$ws.Range("A22").Value = 20
$ws.Range("C22").Value = "NSE:INDIANHUME"
$ws.Range("A23").Value = 21
$ws.Range("C23").Value = "NSE:INDIGO"
$ws.Range("A24").Value = 22
$ws.Range("C24").Value = "NSE:INDORAMA"
$ws.Range("A25").Value = 23
$ws.Range("C25").Value = "NSE:IRMENERGY"
$ws.Range("A26").Value = 24
$ws.Range("C26").Value = "NSE:ITDC"
$ws.Range("A27").Value = 25
$ws.Range("C27").Value = "NSE:JAYSREETEA"
$ws.Range("A28").Value = 26
$ws.Range("C28").Value = "NSE:KAKATCEM"
$ws.Range("A29").Value = 27
$ws.Range("C29").Value = "NSE:KMSUGAR"
$ws.Range("A30").Value = 28
$ws.Range("C30").Value = "NSE:KOTHARIPET"
$ws.Range("A31").Value = 29
$ws.Range("C31").Value = "NSE:KSCL"
$ws.Range("A32").Value = 30
$ws.Range("C32").Value = "NSE:LXCHEM"
$ws.Range("A33").Value = 31
$ws.Range("C33").Value = "NSE:MANGALAM"
$ws.Range("A34").Value = 32
$ws.Range("C34").Value = "NSE:MANGCHEFER"
$ws.Range("A35").Value = 33
$ws.Range("C35").Value = "NSE:MGEL"
$ws.Range("A36").Value = 34
$ws.Range("C36").Value = "NSE:MHLXMIRU"
$ws.Range("A37").Value = 35
$ws.Range("C37").Value = "NSE:MICEL"
$ws.Range("A38").Value = 36
$ws.Range("C38").Value = "NSE:NDTV"
$ws.Range("A39").Value = 37
$ws.Range("C39").Value = "NSE:POWERGRID"
$ws.Range("A40").Value = 38
$ws.Range("C40").Value = "NSE:PRECAM"
$ws.Range("A41").Value = 39
$ws.Range("C41").Value = "NSE:PRINCEPIPE"
$ws.Range("A42").Value = 40
$ws.Range("C42").Value = "NSE:RADIANTCMS"
$ws.Range("A43").Value = 41
$ws.Range("C43").Value = "NSE:RADICO"
$ws.Range("A44").Value = 42
$ws.Range("C44").Value = "NSE:RAMCOCEM"
$ws.Range("A45").Value = 43
$ws.Range("C45").Value = "NSE:RPSGVENT"
$ws.Range("A46").Value = 44
$ws.Range("C46").Value = "NSE:RUPA"
$ws.Range("A47").Value = 45
$ws.Range("C47").Value = "NSE:SALASAR"

# Copy style of column A (row 21, which has the index style) down to the new A cells
$ws.Range("A21").Copy()
$ws.Range("A22:A47").PasteSpecial(-4122)
$excel.CutCopyMode = $false
